$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 53
$ws.Range("C2").Value = "face/face031.jpg"
$ws.Range("D2").Value = "spielen"
$ws.Range("E2").Value = "face"

$ws.Range("B3").Value = 81
$ws.Range("C3").Value = "house/house019.jpg"
$ws.Range("D3").Value = "husten"
$ws.Range("E3").Value = "house"

$ws.Range("B4").Value = 126
$ws.Range("C4").Value = "house/house008.jpg"
$ws.Range("D4").Value = "posten"
$ws.Range("E4").Value = "house"

$ws.Range("B5").Value = 118
$ws.Range("C5").Value = "house/house010.jpg"
$ws.Range("D5").Value = "kehren"
$ws.Range("E5").Value = "house"

$ws.Range("B6").Value = 88
$ws.Range("C6").Value = "face/face019.jpg"
$ws.Range("D6").Value = "scheitern"
$ws.Range("E6").Value = "face"

$ws.Range("B7").Value = 43
$ws.Range("C7").Value = "house/house017.jpg"
$ws.Range("D7").Value = "schätzen"
$ws.Range("E7").Value = "house"

$ws.Range("B8").Value = 55
$ws.Range("C8").Value = "face/face030.jpg"
$ws.Range("D8").Value = "drehen"
$ws.Range("E8").Value = "face"

$ws.Range("B9").Value = 12
$ws.Range("C9").Value = "house/house009.jpg"
$ws.Range("D9").Value = "ehren"
$ws.Range("E9").Value = "house"

$ws.Range("B10").Value = 13
$ws.Range("C10").Value = "face/face022.jpg"
$ws.Range("D10").Value = "fühlen"
$ws.Range("E10").Value = "face"

$ws.Range("B11").Value = 104
$ws.Range("C11").Value = "house/house026.jpg"
$ws.Range("D11").Value = "loben"
$ws.Range("E11").Value = "house"

$ws.Range("B12").Value = 44
$ws.Range("C12").Value = "house/house023.jpg"
$ws.Range("D12").Value = "starten"
$ws.Range("E12").Value = "house"

$ws.Range("B13").Value = 45
$ws.Range("C13").Value = "house/house014.jpg"
$ws.Range("D13").Value = "hupen"
$ws.Range("E13").Value = "house"

$ws.Range("B14").Value = 115
$ws.Range("C14").Value = "face/face000.jpg"
$ws.Range("D14").Value = "kaufen"
$ws.Range("E14").Value = "face"

$ws.Range("B15").Value = 111
$ws.Range("C15").Value = "face/face009.jpg"
$ws.Range("D15").Value = "regnen"
$ws.Range("E15").Value = "face"

$ws.Range("B16").Value = 87
$ws.Range("C16").Value = "face/face002.jpg"
$ws.Range("D16").Value = "haken"
$ws.Range("E16").Value = "face"

$ws.Range("B17").Value = 47
$ws.Range("C17").Value = "house/house030.jpg"
$ws.Range("D17").Value = "wiegen"
$ws.Range("E17").Value = "house"

$ws.Range("B18").Value = 49
$ws.Range("C18").Value = "house/house021.jpg"
$ws.Range("D18").Value = "gelten"
$ws.Range("E18").Value = "house"

$ws.Range("B19").Value = 75
$ws.Range("C19").Value = "house/house002.jpg"
$ws.Range("D19").Value = "schenken"
$ws.Range("E19").Value = "house"

$ws.Range("B20").Value = 72
$ws.Range("C20").Value = "face/face020.jpg"
$ws.Range("D20").Value = "töten"
$ws.Range("E20").Value = "face"

$ws.Range("B21").Value = 84
$ws.Range("C21").Value = "house/house020.jpg"
$ws.Range("D21").Value = "bleiben"
$ws.Range("E21").Value = "house"

$ws.Range("B22").Value = 61
$ws.Range("C22").Value = "house/house012.jpg"
$ws.Range("D22").Value = "bitten"
$ws.Range("E22").Value = "house"

$ws.Range("B23").Value = 119
$ws.Range("C23").Value = "face/face008.jpg"
$ws.Range("D23").Value = "währen"
$ws.Range("E23").Value = "face"

$ws.Range("B24").Value = 85
$ws.Range("C24").Value = "face/face024.jpg"
$ws.Range("D24").Value = "dauern"
$ws.Range("E24").Value = "face"

$ws.Range("B25").Value = 16
$ws.Range("C25").Value = "face/face018.jpg"
$ws.Range("D25").Value = "füllen"
$ws.Range("E25").Value = "face"

$ws.Range("B26").Value = 10
$ws.Range("C26").Value = "house/house029.jpg"
$ws.Range("D26").Value = "laufen"
$ws.Range("E26").Value = "house"

$ws.Range("B27").Value = 102
$ws.Range("C27").Value = "face/face021.jpg"
$ws.Range("D27").Value = "hoffen"
$ws.Range("E27").Value = "face"

$ws.Range("B28").Value = 80
$ws.Range("C28").Value = "house/house024.jpg"
$ws.Range("D28").Value = "krachen"
$ws.Range("E28").Value = "house"

$ws.Range("B29").Value = 14
$ws.Range("C29").Value = "house/house011.jpg"
$ws.Range("D29").Value = "klappen"
$ws.Range("E29").Value = "house"

$ws.Range("B30").Value = 99
$ws.Range("C30").Value = "face/face007.jpg"
$ws.Range("D30").Value = "gründen"
$ws.Range("E30").Value = "face"

$ws.Range("B31").Value = 48
$ws.Range("C31").Value = "face/face026.jpg"
$ws.Range("D31").Value = "sieben"
$ws.Range("E31").Value = "face"

$ws.Range("B32").Value = 120
$ws.Range("C32").Value = "face/face027.jpg"
$ws.Range("D32").Value = "rasen"
$ws.Range("E32").Value = "face"

$ws.Range("B33").Value = 78
$ws.Range("C33").Value = "face/face028.jpg"
$ws.Range("D33").Value = "hauen"
$ws.Range("E33").Value = "face"
